# Clear the "Study not disseminated" evaluation-methods indicator cells
# (G6:G7) on the "Scoring POC_v1" sheet. These had been left set to 1 from
# testing; clearing them resets the scoring example back to a blank state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Scoring POC_v1")
$ws.Range("G6:G7").ClearContents()
$ws.Activate() | Out-Null
$ws.Range("G8").Select() | Out-Null
